$d = $word.ActiveDocument

# --- Change 1: merge ")" and " " math runs into a single run in the
#     first oMathPara (Propiedades de la suma) paragraph. This edits
#     OMML (m:t) content that Find/Replace cannot reach, so the whole
#     paragraph is swapped via InsertXML. ---
$p5 = $d.Paragraphs(5)
$xmlPara5 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body>
<w:p w14:paraId="6BCB04E1" w14:textId="296E98BC" w:rsidR="00527D8A" w:rsidRPr="00527D8A" w:rsidRDefault="00733F00" w:rsidP="00527D8A"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:iCs/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>1</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve">) </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve">Asociatividad: </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F022"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> a, b, c </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F0CE"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> Z, (a + b) + c = a + (b + c). </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:br/></m:r></m:oMath><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:sub></m:sSub><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve">) Conmutatividad: </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F022"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> a, b </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F0CE"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> Z, a + b = b + a. </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:br/></m:r></m:oMath><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>3</m:t></m:r></m:sub></m:sSub><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve">) Existencia de neutro: </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F024"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> 0 </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F0CE"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> Z tal que </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F022"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> a </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F0CE"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> Z, a + 0 = a.</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:br/></m:r></m:oMath><m:oMath><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:softHyphen/></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:softHyphen/></m:r></m:e><m:sub><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>4</m:t></m:r></m:sub></m:sSub><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve">) Existencia de opuesto: Para cada a </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F0CE"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> Z, </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F024"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> (</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F02D"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve">a) </m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F0CE"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> Z tal que a + (</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><w:sym w:font="Symbol" w:char="F02D"/></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve">a) = 0.  </m:t></m:r></m:oMath></m:oMathPara></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$p5.Range.InsertXML($xmlPara5)

# --- Change 2: merge the three runs "...la" + " " + "diferencia..."
#     into a single run. Plain text, so Find/Replace handles it and
#     naturally collapses the run into one. ---
$d.Content.Find.Execute(
  "si a divide a dos enteros, entonces divide tanto a la suma como a la diferencia de dichos enteros.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "si a divide a dos enteros, entonces divide tanto a la suma como a la diferencia de dichos enteros.",
  2)

# --- Change 3: split the "mcm" proofErr/run group (adding
#     gramStart/gramEnd around the opening parenthesis) AND append the
#     five new paragraphs (Propiedad 3.10, 3.11, 3.12, TEOREMA 3.9,
#     TEOREMA 3.10) right after the "TEOREMA 3.8" paragraph. Both are
#     done in one InsertXML call: the target paragraph's corrected
#     content is included first (InsertXML on a collapsed range
#     replaces the preceding paragraph's content with the first
#     supplied paragraph, then inserts the rest after it). ---
$p51 = $d.Paragraphs(51)
$r51 = $p51.Range
$r51.Collapse(0)
$xmlTail = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body>
<w:p w14:paraId="494EF1B1" w14:textId="1FEBBBA8" w:rsidR="009A1A74" w:rsidRDefault="009A1A74" w:rsidP="000A4DAA"><w:r><w:t xml:space="preserve">TEOREMA 3.8 Para a, b </w:t></w:r><w:r><w:sym w:font="Symbol" w:char="F0CE"/></w:r><w:r><w:t xml:space="preserve"> Z + se verifica que </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>mcm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">a, b) </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mcd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(a, b) = a b</w:t></w:r></w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Propiedad 3.10 Si b y c son </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>coprimos</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> con a, entonces b c es </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>coprimo</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> con a.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Propiedad 3.11 Sean a y b </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>coprimos</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> y c un entero tal que </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>c|a</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> y </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>c|b</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, entonces c = 1 o c = </w:t>
      </w:r>
      <w:r>
        <w:sym w:font="Symbol" w:char="F02D"/>
      </w:r>
      <w:r>
        <w:t>1.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Propiedad 3.12 Sea p </w:t>
      </w:r>
      <w:r>
        <w:sym w:font="Symbol" w:char="F0CE"/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Z </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>+ .</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Si a, b </w:t>
      </w:r>
      <w:r>
        <w:sym w:font="Symbol" w:char="F0CE"/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Z, p es primo y p divide al producto ab entonces divide a uno de ellos. En símbolos, si p es primo y </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>p|a</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> b, entonces </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>p|a</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> o </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>p|b</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
        </w:rPr>
      </w:pPr>
      <m:oMathPara>
        <m:oMath>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
            </w:rPr>
            <m:t>TEOREMA 3.9 Si p es primo y p</m:t>
          </m:r>
          <m:d>
            <m:dPr>
              <m:begChr m:val="|"/>
              <m:endChr m:val="|"/>
              <m:ctrlPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                </w:rPr>
              </m:ctrlPr>
            </m:dPr>
            <m:e>
              <m:d>
                <m:dPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:dPr>
                <m:e>
                  <m:r>
                    <m:rPr>
                      <m:sty m:val="p"/>
                    </m:rPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                    </w:rPr>
                    <m:t xml:space="preserve"> </m:t>
                  </m:r>
                  <m:sSub>
                    <m:sSubPr>
                      <m:ctrlPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                        </w:rPr>
                      </m:ctrlPr>
                    </m:sSubPr>
                    <m:e>
                      <m:r>
                        <m:rPr>
                          <m:sty m:val="p"/>
                        </m:rPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                        </w:rPr>
                        <m:t>a</m:t>
                      </m:r>
                    </m:e>
                    <m:sub>
                      <m:r>
                        <m:rPr>
                          <m:sty m:val="p"/>
                        </m:rPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                        </w:rPr>
                        <m:t>1</m:t>
                      </m:r>
                    </m:sub>
                  </m:sSub>
                  <m:sSub>
                    <m:sSubPr>
                      <m:ctrlPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                        </w:rPr>
                      </m:ctrlPr>
                    </m:sSubPr>
                    <m:e>
                      <m:r>
                        <m:rPr>
                          <m:sty m:val="p"/>
                        </m:rPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                        </w:rPr>
                        <m:t>a</m:t>
                      </m:r>
                    </m:e>
                    <m:sub>
                      <m:r>
                        <m:rPr>
                          <m:sty m:val="p"/>
                        </m:rPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                        </w:rPr>
                        <m:t>2</m:t>
                      </m:r>
                    </m:sub>
                  </m:sSub>
                  <m:r>
                    <m:rPr>
                      <m:sty m:val="p"/>
                    </m:rPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                    </w:rPr>
                    <m:t xml:space="preserve"> ... </m:t>
                  </m:r>
                  <m:sSub>
                    <m:sSubPr>
                      <m:ctrlPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                        </w:rPr>
                      </m:ctrlPr>
                    </m:sSubPr>
                    <m:e>
                      <m:r>
                        <m:rPr>
                          <m:sty m:val="p"/>
                        </m:rPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                        </w:rPr>
                        <m:t>a</m:t>
                      </m:r>
                    </m:e>
                    <m:sub>
                      <m:r>
                        <m:rPr>
                          <m:sty m:val="p"/>
                        </m:rPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                        </w:rPr>
                        <m:t>n</m:t>
                      </m:r>
                    </m:sub>
                  </m:sSub>
                  <m:r>
                    <m:rPr>
                      <m:sty m:val="p"/>
                    </m:rPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                    </w:rPr>
                    <m:t xml:space="preserve"> </m:t>
                  </m:r>
                </m:e>
              </m:d>
              <m:r>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                </w:rPr>
                <m:t>entonces p</m:t>
              </m:r>
            </m:e>
          </m:d>
          <m:sSub>
            <m:sSubPr>
              <m:ctrlPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                </w:rPr>
              </m:ctrlPr>
            </m:sSubPr>
            <m:e>
              <m:r>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                </w:rPr>
                <m:t>a</m:t>
              </m:r>
            </m:e>
            <m:sub>
              <m:r>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
                </w:rPr>
                <m:t>i</m:t>
              </m:r>
            </m:sub>
          </m:sSub>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
            </w:rPr>
            <m:t xml:space="preserve"> para algún i con 1 </m:t>
          </m:r>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
            </w:rPr>
            <w:sym w:font="Symbol" w:char="F0A3"/>
          </m:r>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
            </w:rPr>
            <m:t xml:space="preserve"> i </m:t>
          </m:r>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
            </w:rPr>
            <w:sym w:font="Symbol" w:char="F0A3"/>
          </m:r>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
            </w:rPr>
            <m:t xml:space="preserve"> n. </m:t>
          </m:r>
        </m:oMath>
      </m:oMathPara>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">TEOREMA 3.10: Teorema fundamental de la aritmética. Si n es un entero mayor que 1, entonces n es primo o puede escribirse de manera única en la forma n = p1 a 1 p2 a 2 ... </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>pk</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> a k donde p1 </w:t>
      </w:r>
      <w:r>
        <w:sym w:font="Symbol" w:char="F03C"/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> p2 </w:t>
      </w:r>
      <w:r>
        <w:sym w:font="Symbol" w:char="F03C"/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> ... </w:t>
      </w:r>
      <w:r>
        <w:sym w:font="Symbol" w:char="F03C"/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>pk</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> son números primos y los </w:t>
      </w:r>
      <w:r>
        <w:t>exponentes a</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>1 ,</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> a2 , ..., </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ak</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> son enteros positivos.</w:t>
      </w:r>
    </w:p>

</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$r51.InsertXML($xmlTail)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
